$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# H88/I88 were empty inline-string placeholders; clearing makes the cells truly blank (no cell at all).
$ws.Range("H88:I88").ClearContents()

# Append ticket rows 89-112 (new incident log entries dated 2024-05-21).
# Row 89
$ws.Range("A89").NumberFormat = "@"
$ws.Range("A89").Value = '2024-05-21'
$ws.Range("A89").ClearFormats()
$ws.Range("B89").Value = '09:25:35'
$ws.Range("C89").Value = '-'
$ws.Range("D89").Value = 'Detección de sealling mal puesto'
$ws.Range("E89").Value = '-'
$ws.Range("F89").Value = '-'
$ws.Range("G89").Value = '-'
$ws.Range("H89").Value = '09:25:42'
$ws.Range("I89").Value = '0:00:07'

# Row 90
$ws.Range("A90").NumberFormat = "@"
$ws.Range("A90").Value = '2024-05-21'
$ws.Range("A90").ClearFormats()
$ws.Range("B90").Value = '09:25:44'
$ws.Range("C90").Value = '-'
$ws.Range("D90").Value = 'No detecta presencia power CP'
$ws.Range("E90").Value = '-'
$ws.Range("F90").Value = '-'
$ws.Range("G90").Value = '-'

# Row 91
$ws.Range("A91").NumberFormat = "@"
$ws.Range("A91").Value = '2024-05-21'
$ws.Range("A91").ClearFormats()
$ws.Range("B91").Value = '09:26:01'
$ws.Range("C91").Value = '-'
$ws.Range("D91").Value = 'AOI (fallo etiqueta)'
$ws.Range("E91").Value = '-'
$ws.Range("F91").Value = '-'
$ws.Range("G91").Value = '-'
$ws.Range("H91").Value = '09:26:07'
$ws.Range("I91").Value = '0:00:06'

# Row 92
$ws.Range("A92").NumberFormat = "@"
$ws.Range("A92").Value = '2024-05-21'
$ws.Range("A92").ClearFormats()
$ws.Range("B92").Value = '09:28:37'
$ws.Range("C92").Value = 'AOI no detecta pieza'
$ws.Range("D92").Value = '-'
$ws.Range("E92").Value = '-'
$ws.Range("F92").Value = '-'
$ws.Range("G92").Value = '-'

# Row 93
$ws.Range("A93").NumberFormat = "@"
$ws.Range("A93").Value = '2024-05-21'
$ws.Range("A93").ClearFormats()
$ws.Range("B93").Value = '09:28:48'
$ws.Range("C93").Value = 'Ascensor no sube'
$ws.Range("D93").Value = '-'
$ws.Range("E93").Value = '-'
$ws.Range("F93").Value = '-'
$ws.Range("G93").Value = '-'
$ws.Range("H93").Value = '09:28:57'
$ws.Range("I93").Value = '0:00:09'

# Row 94
$ws.Range("A94").NumberFormat = "@"
$ws.Range("A94").Value = '2024-05-21'
$ws.Range("A94").ClearFormats()
$ws.Range("B94").Value = '09:28:59'
$ws.Range("C94").Value = 'No atornilla clips'
$ws.Range("D94").Value = '-'
$ws.Range("E94").Value = '-'
$ws.Range("F94").Value = '-'
$ws.Range("G94").Value = '-'
$ws.Range("H94").Value = '09:29:04'
$ws.Range("I94").Value = '0:00:05'

# Row 95
$ws.Range("A95").NumberFormat = "@"
$ws.Range("A95").Value = '2024-05-21'
$ws.Range("A95").ClearFormats()
$ws.Range("B95").Value = '09:29:32'
$ws.Range("C95").Value = 'No coge placa'
$ws.Range("D95").Value = '-'
$ws.Range("E95").Value = '-'
$ws.Range("F95").Value = '-'
$ws.Range("G95").Value = '-'
$ws.Range("H95").Value = '09:29:37'
$ws.Range("I95").Value = '0:00:05'

# Row 96
$ws.Range("A96").NumberFormat = "@"
$ws.Range("A96").Value = '2024-05-21'
$ws.Range("A96").ClearFormats()
$ws.Range("B96").Value = '09:36:37'
$ws.Range("C96").Value = '-'
$ws.Range("D96").Value = '-'
$ws.Range("E96").Value = '-'
$ws.Range("F96").Value = 'Core enganchado'
$ws.Range("G96").Value = '-'

# Row 97
$ws.Range("A97").NumberFormat = "@"
$ws.Range("A97").Value = '2024-05-21'
$ws.Range("A97").ClearFormats()
$ws.Range("B97").Value = '09:37:53'
$ws.Range("C97").Value = 'Palet atascado en la curva'
$ws.Range("D97").Value = '-'
$ws.Range("E97").Value = '-'
$ws.Range("F97").Value = '-'
$ws.Range("G97").Value = '-'

# Row 98
$ws.Range("A98").NumberFormat = "@"
$ws.Range("A98").Value = '2024-05-21'
$ws.Range("A98").ClearFormats()
$ws.Range("B98").Value = '09:39:42'
$ws.Range("C98").Value = '-'
$ws.Range("D98").Value = 'Cámara no detecta Top cover'
$ws.Range("E98").Value = '-'
$ws.Range("F98").Value = '-'
$ws.Range("G98").Value = '-'

# Row 99
$ws.Range("A99").NumberFormat = "@"
$ws.Range("A99").Value = '2024-05-21'
$ws.Range("A99").ClearFormats()
$ws.Range("B99").Value = '09:42:43'
$ws.Range("C99").Value = 'Fallo en paletizador'
$ws.Range("D99").Value = '-'
$ws.Range("E99").Value = '-'
$ws.Range("F99").Value = '-'
$ws.Range("G99").Value = '-'

# Row 100
$ws.Range("A100").NumberFormat = "@"
$ws.Range("A100").Value = '2024-05-21'
$ws.Range("A100").ClearFormats()
$ws.Range("B100").Value = '09:43:11'
$ws.Range("C100").Value = '-'
$ws.Range("D100").Value = '-'
$ws.Range("E100").Value = 'Power atascado en prensa, cuesta sacar'
$ws.Range("F100").Value = '-'
$ws.Range("G100").Value = '-'

# Row 101
$ws.Range("A101").NumberFormat = "@"
$ws.Range("A101").Value = '2024-05-21'
$ws.Range("A101").ClearFormats()
$ws.Range("B101").Value = '09:49:24'
$ws.Range("C101").Value = 'No atornilla tapa'
$ws.Range("D101").Value = '-'
$ws.Range("E101").Value = '-'
$ws.Range("F101").Value = '-'
$ws.Range("G101").Value = '-'

# Row 102
$ws.Range("A102").NumberFormat = "@"
$ws.Range("A102").Value = '2024-05-21'
$ws.Range("A102").ClearFormats()
$ws.Range("B102").Value = '10:40:56'
$ws.Range("C102").Value = 'Palet atascado en la curva'
$ws.Range("D102").Value = '-'
$ws.Range("E102").Value = '-'
$ws.Range("F102").Value = '-'
$ws.Range("G102").Value = '-'

# Row 103
$ws.Range("A103").NumberFormat = "@"
$ws.Range("A103").Value = '2024-05-21'
$ws.Range("A103").ClearFormats()
$ws.Range("B103").Value = '10:49:33'
$ws.Range("C103").Value = '-'
$ws.Range("D103").Value = 'Cámara no detecta foams'
$ws.Range("E103").Value = '-'
$ws.Range("F103").Value = '-'
$ws.Range("G103").Value = '-'

# Row 104
$ws.Range("A104").NumberFormat = "@"
$ws.Range("A104").Value = '2024-05-21'
$ws.Range("A104").ClearFormats()
$ws.Range("B104").Value = '10:59:41'
$ws.Range("C104").Value = '-'
$ws.Range("D104").Value = 'Cámara no detecta Pcb'
$ws.Range("E104").Value = '-'
$ws.Range("F104").Value = '-'
$ws.Range("G104").Value = '-'
$ws.Range("H104").Value = '10:59:43'
$ws.Range("I104").Value = '0:00:02'

# Row 105
$ws.Range("A105").NumberFormat = "@"
$ws.Range("A105").Value = '2024-05-21'
$ws.Range("A105").ClearFormats()
$ws.Range("B105").Value = '11:00:00'
$ws.Range("C105").Value = '-'
$ws.Range("D105").Value = 'Cámara no detecta Power CP'
$ws.Range("E105").Value = '-'
$ws.Range("F105").Value = '-'
$ws.Range("G105").Value = '-'
$ws.Range("H105").Value = '11:00:02'
$ws.Range("I105").Value = '0:00:02'

# Row 106
$ws.Range("A106").NumberFormat = "@"
$ws.Range("A106").Value = '2024-05-21'
$ws.Range("A106").ClearFormats()
$ws.Range("B106").Value = '11:00:10'
$ws.Range("C106").Value = '-'
$ws.Range("D106").Value = 'Cámara no detecta foam derecho'
$ws.Range("E106").Value = '-'
$ws.Range("F106").Value = '-'
$ws.Range("G106").Value = '-'

# Row 107
$ws.Range("A107").NumberFormat = "@"
$ws.Range("A107").Value = '2024-05-21'
$ws.Range("A107").ClearFormats()
$ws.Range("B107").Value = '11:00:12'
$ws.Range("C107").Value = '-'
$ws.Range("D107").Value = 'Power atascado en prensa, cuesta sacar'
$ws.Range("E107").Value = '-'
$ws.Range("F107").Value = '-'
$ws.Range("G107").Value = '-'
$ws.Range("H107").Value = '11:00:14'
$ws.Range("I107").Value = '0:00:02'

# Row 108
$ws.Range("A108").NumberFormat = "@"
$ws.Range("A108").Value = '2024-05-21'
$ws.Range("A108").ClearFormats()
$ws.Range("B108").Value = '11:00:43'
$ws.Range("C108").Value = '-'
$ws.Range("D108").Value = 'Tornillo atascado en tolva'
$ws.Range("E108").Value = '-'
$ws.Range("F108").Value = '-'
$ws.Range("G108").Value = '-'
$ws.Range("H108").Value = '11:00:45'
$ws.Range("I108").Value = '0:00:02'

# Row 109
$ws.Range("A109").NumberFormat = "@"
$ws.Range("A109").Value = '2024-05-21'
$ws.Range("A109").ClearFormats()
$ws.Range("B109").Value = '11:02:27'
$ws.Range("C109").Value = 'Etiquetadora'
$ws.Range("D109").Value = '-'
$ws.Range("E109").Value = '-'
$ws.Range("F109").Value = '-'
$ws.Range("G109").Value = '-'
$ws.Range("H109").Value = '11:02:28'
$ws.Range("I109").Value = '0:00:01'

# Row 110
$ws.Range("A110").NumberFormat = "@"
$ws.Range("A110").Value = '2024-05-21'
$ws.Range("A110").ClearFormats()
$ws.Range("B110").Value = '11:06:41'
$ws.Range("C110").Value = '-'
$ws.Range("D110").Value = '-'
$ws.Range("E110").Value = '-'
$ws.Range("F110").Value = 'Cover atascado'
$ws.Range("G110").Value = '-'
$ws.Range("H110").Value = '11:06:48'
$ws.Range("I110").Value = '0:00:07'

# Row 111
$ws.Range("A111").NumberFormat = "@"
$ws.Range("A111").Value = '2024-05-21'
$ws.Range("A111").ClearFormats()
$ws.Range("B111").Value = '11:09:15'
$ws.Range("C111").Value = 'Fallo en paletizador'
$ws.Range("D111").Value = '-'
$ws.Range("E111").Value = '-'
$ws.Range("F111").Value = '-'
$ws.Range("G111").Value = '-'
$ws.Range("H111").Value = '11:09:17'
$ws.Range("I111").Value = '0:00:02'

# Row 112
$ws.Range("A112").NumberFormat = "@"
$ws.Range("A112").Value = '2024-05-21'
$ws.Range("A112").ClearFormats()
$ws.Range("B112").Value = '11:10:39'
$ws.Range("C112").Value = 'Fallo fijador tapa'
$ws.Range("D112").Value = '-'
$ws.Range("E112").Value = '-'
$ws.Range("F112").Value = '-'
$ws.Range("G112").Value = '-'
$ws.Range("H112").Value = '11:10:40'
$ws.Range("I112").Value = '0:00:01'

